# Re-upload edit: correct the "Ujaran Minta-Minta" / "Komentar Lain" counts
# and percentages for data row index 8/9 (A10/A11), matching the corrected
# source numbers (423 + 7991 = 8414 total), and refresh the sheet
# selection/scroll state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 10 (Kategori index 8): count 680 -> 423, percentage recalculated ---
$ws.Range("C10").Value = 423
$ws.Range("D10").Value = 5.02733

# --- Row 11 (Kategori index 9): count 7732 -> 7991, percentage recalculated ---
$ws.Range("C11").Value = 7991
$ws.Range("D11").Value = 94.972665

# The percentage cell for the corrected row picked up a wrap-text format
# in the saved workbook.
$ws.Range("D10").WrapText = $true

# Refresh the view: scroll back to the top of the sheet and leave the
# selection on the last-edited cell (E11).
[void]$ws.Range("A1").Select()
[void]$ws.Range("E11").Select()
